$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-09-08T15:07:54"
$ws.Range("S4").Value = 36.22
$ws.Range("T4").Value = 46.26
$ws.Range("U4").Value = 37
$ws.Range("V4").Value = 36.74
$ws.Range("W4").Value = 36.69
$ws.Range("S6").Value = -1.56
$ws.Range("T6").Value = -1.85
$ws.Range("U6").Value = -1.52
$ws.Range("V6").Value = -1.47
$ws.Range("W6").Value = -1.28
$ws.Range("X6").Value = -1.36
$ws.Range("T9").Value = 47.17
$ws.Range("U9").Value = 37.29
$ws.Range("V9").Value = 37.42
$ws.Range("W9").Value = 37.6
$ws.Range("S11").Value = -0.85
$ws.Range("T11").Value = -0.9399999999999999
$ws.Range("U11").Value = -1.23
$ws.Range("V11").Value = -0.79
$ws.Range("W11").Value = -0.38
$ws.Range("X11").Value = 0
$ws.Range("S12").Value = 4.07
$ws.Range("T14").Value = 47.17
$ws.Range("U14").Value = 37.29
$ws.Range("V14").Value = 37.46
$ws.Range("W14").Value = 37.6
$ws.Range("Y14").Value = 34.9
$ws.Range("S16").Value = -0.85
$ws.Range("T16").Value = -0.9399999999999999
$ws.Range("U16").Value = -1.23
$ws.Range("V16").Value = -0.75
$ws.Range("W16").Value = -0.38
$ws.Range("X16").Value = 0
$ws.Range("Y16").Value = 0.1
$ws.Range("S17").Value = 4.07
$ws.Range("S19").Value = 36.01
$ws.Range("T19").Value = 46
$ws.Range("U19").Value = 36.79
$ws.Range("V19").Value = 36.56
$ws.Range("W19").Value = 36.59
$ws.Range("X19").Value = 36.55
$ws.Range("Y19").Value = 33.58
$ws.Range("T20").Value = 0
$ws.Range("S21").Value = -1.76
$ws.Range("T21").Value = -2.12
$ws.Range("U21").Value = -1.73
$ws.Range("V21").Value = -1.65
$ws.Range("W21").Value = -1.39
$ws.Range("X21").Value = -1.43
$ws.Range("Y21").Value = -1.21
$ws.Range("S24").Value = 36.01
$ws.Range("T24").Value = 46
$ws.Range("U24").Value = 36.79
$ws.Range("V24").Value = 36.56
$ws.Range("W24").Value = 36.59
$ws.Range("X24").Value = 36.55
$ws.Range("Y24").Value = 33.58
$ws.Range("S26").Value = -1.76
$ws.Range("T26").Value = -2.12
$ws.Range("U26").Value = -1.73
$ws.Range("V26").Value = -1.65
$ws.Range("W26").Value = -1.39
$ws.Range("X26").Value = -1.43
$ws.Range("Y26").Value = -1.21
$ws.Range("S29").Value = 35.77
$ws.Range("T29").Value = 45.65
$ws.Range("U29").Value = 36.58
$ws.Range("V29").Value = 36.32
$ws.Range("W29").Value = 36.41
$ws.Range("X29").Value = 36.41
$ws.Range("Y29").Value = 33.49
$ws.Range("Z29").Value = 34.63
$ws.Range("T30").Value = 0
$ws.Range("S31").Value = -2
$ws.Range("T31").Value = -2.47
$ws.Range("U31").Value = -1.94
$ws.Range("V31").Value = -1.89
$ws.Range("W31").Value = -1.57
$ws.Range("X31").Value = -1.57
$ws.Range("Y31").Value = -1.31
$ws.Range("Z31").Value = -1
$ws.Range("S34").Value = 37.67
$ws.Range("T34").Value = 61
$ws.Range("U34").Value = 37.77
$ws.Range("V34").Value = 38.09
$ws.Range("W34").Value = 38.32
$ws.Range("X34").Value = 38.91
$ws.Range("Z34").Value = 37
$ws.Range("S36").Value = -0.37
$ws.Range("T36").Value = -0.14
$ws.Range("U36").Value = -0.76
$ws.Range("V36").Value = -0.11
$ws.Range("W36").Value = 0.34
$ws.Range("X36").Value = 0.93
$ws.Range("Z36").Value = 1.37
$ws.Range("S37").Value = 0.27
$ws.Range("T37").Value = 13.03
$ws.Range("S39").Value = 36.22
$ws.Range("T39").Value = 46.26
$ws.Range("U39").Value = 37
$ws.Range("V39").Value = 36.74
$ws.Range("W39").Value = 36.69
$ws.Range("S41").Value = -1.56
$ws.Range("T41").Value = -1.85
$ws.Range("U41").Value = -1.52
$ws.Range("V41").Value = -1.47
$ws.Range("W41").Value = -1.28
$ws.Range("X41").Value = -1.36
$ws.Range("S44").Value = 37.63
$ws.Range("T44").Value = 48.41
$ws.Range("U44").Value = 38.56
$ws.Range("V44").Value = 38.24
$ws.Range("W44").Value = 38.21
$ws.Range("X44").Value = 38.02
$ws.Range("S46").Value = -0.15
$ws.Range("T46").Value = 0.29
$ws.Range("W46").Value = 0.23
$ws.Range("S49").Value = 39.56
$ws.Range("T49").Value = 47.45
$ws.Range("U49").Value = 39.92
$ws.Range("V49").Value = 39.47
$ws.Range("W49").Value = 39.03
$ws.Range("X49").Value = 39.28
$ws.Range("Z49").Value = 36.81
$ws.Range("S51").Value = 1.78
$ws.Range("T51").Value = -0.66
$ws.Range("U51").Value = 1.4
$ws.Range("V51").Value = 1.26
$ws.Range("W51").Value = 1.05
$ws.Range("X51").Value = 1.3
$ws.Range("Z51").Value = 1.18
$ws.Range("S54").Value = 37.97
$ws.Range("T54").Value = 47.88
$ws.Range("U54").Value = 38.44
$ws.Range("V54").Value = 38.79
$ws.Range("W54").Value = 38.56
$ws.Range("X54").Value = 38.4
$ws.Range("S56").Value = 0.19
$ws.Range("T56").Value = -0.24
$ws.Range("U56").Value = -0.08
$ws.Range("V56").Value = 0.58
$ws.Range("W56").Value = 0.58
$ws.Range("S59").Value = 38.91
$ws.Range("T59").Value = 50.07
$ws.Range("U59").Value = 39.75
$ws.Range("V59").Value = 39.43
$ws.Range("W59").Value = 39.44
$ws.Range("X59").Value = 39.23
$ws.Range("S61").Value = 1.13
$ws.Range("T61").Value = 1.95
$ws.Range("U61").Value = 1.23
$ws.Range("V61").Value = 1.22
$ws.Range("W61").Value = 1.46
$ws.Range("X61").Value = 1.26
$ws.Range("S64").Value = 39.68
$ws.Range("T64").Value = 50.97
$ws.Range("U64").Value = 40.51
$ws.Range("V64").Value = 40.13
$ws.Range("W64").Value = 40.1
$ws.Range("X64").Value = 39.89
$ws.Range("S66").Value = 1.9
$ws.Range("T66").Value = 2.85
$ws.Range("U66").Value = 1.98
$ws.Range("V66").Value = 1.93
$ws.Range("W66").Value = 2.13
$ws.Range("S69").Value = 40.15
$ws.Range("T69").Value = 51.46
$ws.Range("U69").Value = 40.68
$ws.Range("V69").Value = 40.34
$ws.Range("W69").Value = 40.53
$ws.Range("X69").Value = 40.36
$ws.Range("S71").Value = 2.37
$ws.Range("T71").Value = 3.34
$ws.Range("U71").Value = 2.16
$ws.Range("V71").Value = 2.14
$ws.Range("W71").Value = 2.55
$ws.Range("S74").Value = 39.11
$ws.Range("T74").Value = 50.02
$ws.Range("U74").Value = 39.84
$ws.Range("V74").Value = 39.39
$ws.Range("W74").Value = 39.4
$ws.Range("X74").Value = 39.19
$ws.Range("Y74").Value = 35.98
$ws.Range("S76").Value = 1.33
$ws.Range("T76").Value = 1.9
$ws.Range("U76").Value = 1.31
$ws.Range("V76").Value = 1.18
$ws.Range("W76").Value = 1.42
$ws.Range("X76").Value = 1.22
$ws.Range("Y76").Value = 1.19
$ws.Range("S79").Value = 39.33
$ws.Range("T79").Value = 50.32
$ws.Range("U79").Value = 39.99
$ws.Range("V79").Value = 39.63
$ws.Range("W79").Value = 39.66
$ws.Range("X79").Value = 39.44
$ws.Range("S81").Value = 1.56
$ws.Range("T81").Value = 2.21
$ws.Range("U81").Value = 1.47
$ws.Range("V81").Value = 1.42
$ws.Range("W81").Value = 1.68
$ws.Range("S84").Value = 38.04
$ws.Range("T84").Value = 46.99
$ws.Range("U84").Value = 37.84
$ws.Range("V84").Value = 39.14
$ws.Range("W84").Value = 39.03
$ws.Range("X84").Value = 38.75
$ws.Range("Z84").Value = 36.43
$ws.Range("S86").Value = 0.27
$ws.Range("T86").Value = -1.13
$ws.Range("U86").Value = -0.68
$ws.Range("V86").Value = 0.9399999999999999
$ws.Range("W86").Value = 1.05
$ws.Range("X86").Value = 0.78
$ws.Range("Z86").Value = 0.8
$ws.Range("S89").Value = 35.77
$ws.Range("T89").Value = 45.65
$ws.Range("U89").Value = 36.58
$ws.Range("V89").Value = 36.32
$ws.Range("W89").Value = 36.41
$ws.Range("X89").Value = 36.41
$ws.Range("Y89").Value = 33.49
$ws.Range("Z89").Value = 34.63
$ws.Range("T90").Value = 0
$ws.Range("S91").Value = -2
$ws.Range("T91").Value = -2.47
$ws.Range("U91").Value = -1.94
$ws.Range("V91").Value = -1.89
$ws.Range("W91").Value = -1.57
$ws.Range("X91").Value = -1.57
$ws.Range("Y91").Value = -1.31
$ws.Range("Z91").Value = -1
